$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Colors (OLE BGR-packed RGB integers)
# -----------------------------------------------------------------
$yellow = 65535      # FFFF00
$green  = 5296274    # 92D050

# -----------------------------------------------------------------
# Row 1 (new row): A1 = "Game Start", bold + orange(theme) fill, no alignment
# Copy format from O4 which already carries that exact style (bold, fillId
# theme9, no alignment) so we reuse the existing themed fill instead of
# creating a new RGB-based one.
# -----------------------------------------------------------------
$ws.Range("A1").Value2 = "Game Start"
$ws.Range("O4").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# -----------------------------------------------------------------
# Row 2: A2 total health, B2 = 10000
# -----------------------------------------------------------------
$ws.Range("A2").Value2 = "Total health"
$ws.Range("B2").Value2 = 10000

# -----------------------------------------------------------------
# Row 3: A3 hunger (trailing space kept), B3 = 5000
# -----------------------------------------------------------------
$ws.Range("A3").Value2 = "Hunger "
$ws.Range("B3").Value2 = 5000

# -----------------------------------------------------------------
# Row 4: A4 coins, B4 = 200 ; I4/J4/K4 new headers (bold+yellow+wrap)
# -----------------------------------------------------------------
$ws.Range("A4").Value2 = "Coins"
$ws.Range("B4").Value2 = 200

$ws.Range("I4").Value2 = "Hunger recovery"
$ws.Range("I4").Font.Bold = $true
$ws.Range("I4").Interior.Color = $yellow
$ws.Range("I4").WrapText = $true

$ws.Range("J4").Value2 = "Heath recovered"
$ws.Range("J4").Font.Bold = $true
$ws.Range("J4").Interior.Color = $yellow
$ws.Range("J4").WrapText = $true

$ws.Range("K4").Value2 = "Health depletion"
$ws.Range("K4").Font.Bold = $true
$ws.Range("K4").Interior.Color = $yellow
$ws.Range("K4").WrapText = $true

# -----------------------------------------------------------------
# Row 5: new data point I5 = 20, shift old I5(5) into J5
# -----------------------------------------------------------------
$ws.Range("I5").Value2 = 20
$ws.Range("J5").Value2 = 5

# -----------------------------------------------------------------
# Row 6: B6 becomes "ENEMIES" (style already bold/green/center - untouched)
#        new data point I6 = 22, shift old I6(7) into J6
# -----------------------------------------------------------------
$ws.Range("B6").Value2 = "ENEMIES"
$ws.Range("I6").Value2 = 22
$ws.Range("J6").Value2 = 7

# -----------------------------------------------------------------
# Row 7: row height becomes 18 (custom) ; B7/C7/D7 new green headers
#        new data point I7 = 30, shift old I7(12) into J7
# -----------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 18

$ws.Range("B7").Value2 = "Health"
$ws.Range("B7").Font.Bold = $true
$ws.Range("B7").Interior.Color = $green

$ws.Range("C7").Value2 = "Coins"
$ws.Range("C7").Font.Bold = $true
$ws.Range("C7").Interior.Color = $green

$ws.Range("D7").Value2 = "Damage"
$ws.Range("D7").Font.Bold = $true
$ws.Range("D7").Interior.Color = $green
$ws.Range("D7").WrapText = $true

$ws.Range("I7").Value2 = 30
$ws.Range("J7").Value2 = 12

# -----------------------------------------------------------------
# Row 8: enemy data moves from B.. to A.. with new B/C/D values
#        new data point I8 = 45, shift old I8(30) into J8
# -----------------------------------------------------------------
$ws.Range("A8").Value2 = "Small enemy"
$ws.Range("B8").Value2 = 3000
$ws.Range("C8").Value2 = 50
$ws.Range("D8").Value2 = 20
$ws.Range("I8").Value2 = 45
$ws.Range("J8").Value2 = 30

# -----------------------------------------------------------------
# Row 9: enemy data moves from B.. to A.. with new B/C/D values
#        new data point I9 = 35, shift old I9(25)->J9, old J9(-1)->K9
# -----------------------------------------------------------------
$ws.Range("A9").Value2 = "Medium enemy"
$ws.Range("B9").Value2 = 5000
$ws.Range("C9").Value2 = 80
$ws.Range("D9").Value2 = 50
$ws.Range("I9").Value2 = 35
$ws.Range("J9").Value2 = 25
$ws.Range("K9").Value2 = -1

# -----------------------------------------------------------------
# Row 10: enemy data moves from B.. to A.. with new B/C/D values
# -----------------------------------------------------------------
$ws.Range("A10").Value2 = "Large enemy"
$ws.Range("B10").Value2 = 7000
$ws.Range("C10").Value2 = 100
$ws.Range("D10").Value2 = 100

# -----------------------------------------------------------------
# New rows 15/16: Notes section
# -----------------------------------------------------------------
$ws.Range("A15").Value2 = "Notes:"

$ws.Rows.Item(16).RowHeight = 43.2
$ws.Range("A16").Value2 = "Rotten Meat "
$ws.Range("B16").Value2 = "( -1 every 1 secs after the initial boost)"
$ws.Range("B16").WrapText = $true

# -----------------------------------------------------------------
# New column J width
# -----------------------------------------------------------------
$ws.Columns.Item(10).ColumnWidth = 9.75

# -----------------------------------------------------------------
# Selection matches target view state
# -----------------------------------------------------------------
$ws.Range("E16").Select() | Out-Null

Write-Host "Shop sheet updated"
